# Generate Report for Handoff
# Update the "Latest Handoff"/"Latest Handback" timestamps for the
# 98a512ad-5466-4600-97cb-d27eadf23803 file row (row 6) on each sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-29-19 12:29:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-19 12:29:46"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-19 12:29:49"
